# Insert a new weekly data row for "Fruta, Feria Lagunitas de Puerto Montt - Arándano (blue)"
# at row 34, pushing the existing rows 34-49 down to 35-50.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 34; this shifts rows 34:49 -> 35:50
# (and the sheet's used range / dimension) automatically, carrying formatting along.
$ws.Rows.Item(34).Insert()

# Populate the newly inserted row 34 with the new record's data.
$ws.Range("A34").Value = 4
$ws.Range("B34").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C34").Value = "Los Lagos"
$ws.Range("D34").Value = 44957
$ws.Range("E34").Value = 10
$ws.Range("F34").Value = "Fruta"
$ws.Range("G34").Value = 100101
$ws.Range("H34").Value = "Berries"
$ws.Range("I34").Value = 100101001
$ws.Range("J34").Value = "Arándano (blue)"
$ws.Range("K34").Value = "Sin especificar"
$ws.Range("L34").Value = "Primera"
$ws.Range("M34").Value = 200
$ws.Range("N34").Value = 2000
$ws.Range("O34").Value = 2200
$ws.Range("P34").Value = 2100
$ws.Range("Q34").Value = "$/bandeja 2 kilos"
$ws.Range("R34").Value = "Provincia de Curicó"
$ws.Range("S34").Value = 1050
$ws.Range("T34").Value = 2
